$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatosGenerales")
$ws.Activate()

$ws.Range("B6").Value = "25"
$ws.Range("B6").Select()
